# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.221.03"
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").Value = "3.342.80"
$ws.Range("E3").Value = "  -5.27%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'553.50"
$ws.Range("E5").Value = "  -4.22%  "
$ws.Range("D6").Value = "'173.89"
$ws.Range("E6").Value = "  -3.27%  "
$ws.Range("D7").Value = "'0.614"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.336.09"
$ws.Range("E8").Value = "  -5.30%  "
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'0.624"
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").Value = "'0.162"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "'53.65"
$ws.Range("E12").Value = "  -3.81%  "
$ws.Range("D13").Value = "'0.0000274"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").Value = "'9.03"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "3.862.42"
$ws.Range("E15").Value = "  -5.64%  "
$ws.Range("D16").Value = "'18.38"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("D18").Value = "3.328.16"
$ws.Range("E18").Value = "  -5.56%  "
$ws.Range("D19").Value = "'11.78"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").Value = "64.056.19"
$ws.Range("E20").Value = "  -3.57%  "
$ws.Range("D21").Value = "'0.976"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").Value = "'426.89"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").Value = "'4.83"
$ws.Range("E23").Value = "  +12.17%  "
$ws.Range("D24").Value = "'4.08"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").Value = "'84.07"
$ws.Range("E25").Value = "  -2.16%  "
$ws.Range("D26").Value = "'13.18"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").Value = "'10.67"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("D28").Value = "'2.82"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").Value = "'8.63"
$ws.Range("E29").Value = "  -5.48%  "
$ws.Range("D30").Value = "'29.71"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("D31").Value = "'6.70"
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").Value = "'595.30"
$ws.Range("E32").Value = "  -5.41%  "
$ws.Range("D33").Value = "'11.41"
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("D35").Value = "'58.22"
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -10.18%  "
$ws.Range("D38").Value = "'3.55"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0758"
$ws.Range("E39").Value = "  -6.95%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'35.50"
$ws.Range("E40").Value = "  -5.15%  "
$ws.Range("D41").Value = "'0.365"
$ws.Range("E41").Value = "  -5.25%  "
$ws.Range("D42").Value = "3.095.73"
$ws.Range("E42").Value = "  -4.98%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("D44").Value = "'2.80"
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("D45").Value = "'3.22"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("D46").Value = "'0.0406"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("D47").Value = "'2.45"
$ws.Range("E47").Value = "  -4.39%  "
$ws.Range("D48").Value = "'0.129"
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("E49").Value = "  -4.55%  "
$ws.Range("D50").Value = "'134.16"
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("D51").Value = "'8.21"
$ws.Range("E51").Value = "  -5.36%  "
